$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enter the new cell text values first, in the same order the original
#     author typed them (this governs shared-string table ordering) ---
$ws.Range("C27").Value = "25 First Class Seats "
$ws.Range("D27").Value = "125 Coach Seats "
$ws.Range("D28").Value = "C >=40"
$ws.Range("B28").Value = "For C1, C2 and C3 "
$ws.Range("E28").Value = "FC + C <=150"
$ws.Range("B27").Value = "The Constraints are met"
$ws.Range("E27").Value = "Total seatls 25+125 = 150 "
$ws.Range("C28").Value = "FC>=25"

# --- Row 27 formatting: centered text, green highlight on the figures ---
$ws.Range("B27:E27").HorizontalAlignment = -4108
$ws.Range("C27:E27").Interior.Color = 5296274

# --- Row 28 formatting: centered text, no fill ---
$ws.Range("B28:E28").HorizontalAlignment = -4108

# --- Update view / selection to match the saved state ---
$ws.Activate()
$ws.Range("C29").Select()
